# Generate Report for Archive
# Updates the localization status from "Ready for handoff" to "In Translation"
# across the Overview / zh-cn / de-de sheets, and lets the status column
# shrink to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) and de-de (col F) status columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.43
$overview.Columns.Item(6).ColumnWidth = 12.43

# --- zh-cn sheet: Status column (col C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.43

# --- de-de sheet: Status column (col C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.43
